$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.400.97"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.574.79"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.23"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3773"
$ws.Range("E7").Value = "  +2.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.98"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3425"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.168"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07690"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.39"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.995"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.938"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.573.92"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.63"
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06739"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.246"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5280"
$ws.Range("E23").Value = "  -3.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.03"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.393.25"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.395"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.785"
$ws.Range("E27").Value = "  -5.15%  "
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.02"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.080"
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.45"
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.747.18"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.248"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.018"
$ws.Range("E34").Value = "  +5.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.028"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08533"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02564"
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2327"
$ws.Range("E39").Value = "  +1.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06578"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.507"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.300"
$ws.Range("E42").Value = "  +3.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.68"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6467"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.11"
$ws.Range("E45").Value = "  -3.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6043"
$ws.Range("E47").Value = "  +1.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.785"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.305"
$ws.Range("E49").Value = "  +10.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.106"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.60"
$ws.Range("E51").Value = "  +3.28%  "
